$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "View" field header to "Cache"
$ws.Range("F1").Value = "Cache"

# Default the boolean field values in column F (rows 2-25) to FALSE
$ws.Range("F2:F25").Value = $false
